$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(4, 6).Value = 46
$ws.Cells.Item(5, 6).Value = 191
$ws.Cells.Item(6, 6).Value = 3824
$ws.Cells.Item(7, 6).Value = 189
$ws.Cells.Item(10, 6).Value = 81
$ws.Cells.Item(12, 6).Value = 678
$ws.Cells.Item(13, 6).Value = 172
$ws.Cells.Item(14, 6).Value = 941
$ws.Cells.Item(15, 6).Value = 69
$ws.Cells.Item(16, 6).Value = 234
$ws.Cells.Item(19, 6).Value = 104
$ws.Cells.Item(20, 6).Value = 90
$ws.Cells.Item(21, 6).Value = 3383
$ws.Cells.Item(22, 6).Value = 5714
$ws.Cells.Item(24, 6).Value = 24
$ws.Cells.Item(26, 6).Value = 515
$ws.Cells.Item(28, 6).Value = 3330
$ws.Cells.Item(29, 6).Value = 351
$ws.Cells.Item(30, 6).Value = 20
$ws.Cells.Item(31, 6).Value = 2438
$ws.Cells.Item(34, 6).Value = 116
$ws.Cells.Item(35, 6).Value = 198
$ws.Cells.Item(36, 6).Value = 255
$ws.Cells.Item(38, 6).Value = 115
$ws.Cells.Item(39, 6).Value = 1004
$ws.Cells.Item(43, 6).Value = 37
$ws.Cells.Item(44, 6).Value = 462
$ws.Cells.Item(45, 6).Value = 60

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 92

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(4, 6).Value = 46
$ws.Cells.Item(5, 6).Value = 191
$ws.Cells.Item(6, 6).Value = 3824
$ws.Cells.Item(7, 6).Value = 189
$ws.Cells.Item(10, 6).Value = 92
$ws.Cells.Item(11, 6).Value = 81
$ws.Cells.Item(13, 6).Value = 678
$ws.Cells.Item(14, 6).Value = 172
$ws.Cells.Item(15, 6).Value = 941
$ws.Cells.Item(16, 6).Value = 69
$ws.Cells.Item(17, 6).Value = 234
$ws.Cells.Item(20, 6).Value = 104
$ws.Cells.Item(21, 6).Value = 90
$ws.Cells.Item(22, 6).Value = 3383
$ws.Cells.Item(23, 6).Value = 5714
$ws.Cells.Item(25, 6).Value = 24
$ws.Cells.Item(27, 6).Value = 515
$ws.Cells.Item(29, 6).Value = 3330
$ws.Cells.Item(30, 6).Value = 351
$ws.Cells.Item(31, 6).Value = 20
$ws.Cells.Item(32, 6).Value = 2438
$ws.Cells.Item(35, 6).Value = 116
$ws.Cells.Item(36, 6).Value = 198
$ws.Cells.Item(37, 6).Value = 255
$ws.Cells.Item(39, 6).Value = 115
$ws.Cells.Item(40, 6).Value = 1004
$ws.Cells.Item(44, 6).Value = 37
$ws.Cells.Item(45, 6).Value = 462
$ws.Cells.Item(46, 6).Value = 60
